$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the "Hug_Date_Derniere_Soumission_C.WORK" row (row 23) ---
$ws.Range("A23").EntireRow.Delete()

# --- 2) Delete the "physicalb" and "occupation" rows (now rows 31 and 32) ---
$ws.Range("A31:A32").EntireRow.Delete()

# --- 3) Insert 3 new rows before "serocov_work.inc" (now row 23) ---
$ws.Range("A23:A25").EntireRow.Insert()

# --- 4) Fill in the 3 new rows ---
$ws.Range("A23").Value = "date_inclusion"
$ws.Range("B23").Value = "date they filled the Specchio inclusion questionnaire"

$ws.Range("A24").Value = "date_last_submission"
$ws.Range("B24").Value = "the date of their most recent questionnaire submission"

$ws.Range("A25").Value = "beyond_inclusion"
$ws.Range("B25").Value = "did they fill any questionnaire beyond the inclusion questionnaire? TRUE / FALSE"

# match style (wrap text, like the rest of column B) for the new B cells
$ws.Range("B23:B25").Style = $ws.Range("B22").Style

# --- 5) Update "occupational_grouping" description (row 34) ---
$ws.Range("B34").Value = "grouping created by Berg et al. based on 2-digit ISCO-08"

# --- 6) Update "key_occupation" literature (row 35) ---
$ws.Range("C35").Value = "Berg at al. 2023"
$ws.Range("C35").Style = $ws.Range("B35").Style

# --- 7) Update "health_workers" description + literature (row 36), add hyperlink ---
$ws.Range("B36").Value = "definition from WHO classification of ""health workers"" that uses the ISCO-08 4-digit codes"
$ws.Range("C36").Value = "https://www.who.int/publications/m/item/classifying-health-workers"
$ws.Hyperlinks.Add($ws.Range("C36"), "https://www.who.int/publications/m/item/classifying-health-workers") | Out-Null

# --- 8) Turn the A1:C36 range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:C36"), $null, 1)
$tbl.Name = "Tableau2"
$tbl.TableStyle = "TableStyleMedium1"

# --- 9) Update the sheet view selection / scroll position ---
$ws.Range("E34").Select()
$excel.ActiveWindow.ScrollRow = 1
